$d = $word.ActiveDocument

# Paragraphs to append at the end of the document, in order. An empty
# string denotes a blank paragraph.
$newParagraphs = @(
    "APIs ultimately transmit data to and from applications.. ",
    "APIs receive a request from a client, respond with a response obj.",
    "Request objs contain action verbs(POST, GET) and may contain some sort of  content",
    "Response objs instead have a Status Code, and may also have content..",
    "",
    "We start by creating our APIController.. this is where we write out our API operations, our endpoints/methods.. ",
    "We need to create our Food Model class as well"
)

foreach ($text in $newParagraphs) {
    if ($text -eq "") {
        $d.Content.InsertParagraphAfter() | Out-Null
    }
    else {
        $r = $d.Content
        $r.Collapse(0)
        $r.InsertAfter("`r" + $text)
    }
}
